$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, shifting existing rows 150-154 down to 151-155
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with the new weekly price record
$ws.Range("A150").Value = 5
$ws.Range("B150").Value = "Macroferia Regional de Talca"
$ws.Range("C150").Value = "Maule"
$ws.Range("D150").Value = 45239
$ws.Range("E150").Value = 7
$ws.Range("F150").Value = 100112022
$ws.Range("G150").Value = "Arveja Verde"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 400
$ws.Range("K150").Value = 20000
$ws.Range("L150").Value = 20000
$ws.Range("M150").Value = 20000
$ws.Range("N150").Value = '$/saco 25 kilos'
$ws.Range("O150").Value = "Región del Maule"
$ws.Range("P150").Value = 800
$ws.Range("Q150").Value = 25
$ws.Range("R150").Value = "Hortaliza"
